# Update cryptocurrency Price (D) and Volume(1h) (E) columns with refreshed values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.718.28"
$ws.Range("E2").Value = "  -1.59%  "

$ws.Range("D3").Value = "'2.677.19"
$ws.Range("E3").Value = "  -1.31%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'598.81"
$ws.Range("E5").Value = "  -0.22%  "

$ws.Range("D6").Value = "'166.42"
$ws.Range("E6").Value = "  +2.27%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("E8").Value = "  +0.40%  "

$ws.Range("D9").Value = "'2.677.01"
$ws.Range("E9").Value = "  -1.21%  "

$ws.Range("E10").Value = "  +2.36%  "

$ws.Range("E11").Value = "  +1.09%  "

$ws.Range("E12").Value = "  +0.10%  "

$ws.Range("E13").Value = "  -1.66%  "

$ws.Range("D14").Value = "'27.91"
$ws.Range("E14").Value = "  -2.14%  "

$ws.Range("D15").Value = "'3.166.17"
$ws.Range("E15").Value = "  -0.87%  "

$ws.Range("E16").Value = "  -1.73%  "

$ws.Range("D17").Value = "'67.711.43"
$ws.Range("E17").Value = "  -1.40%  "

$ws.Range("D18").Value = "'2.679.94"
$ws.Range("E18").Value = "  -0.39%  "

$ws.Range("D19").Value = "'11.78"
$ws.Range("E19").Value = "  -0.76%  "

$ws.Range("D20").Value = "'7.75"

$ws.Range("D21").Value = "'364.41"
$ws.Range("E21").Value = "  -0.40%  "

$ws.Range("E22").Value = "  -3.29%  "

$ws.Range("E23").Value = "  -2.17%  "

$ws.Range("E24").Value = "  -3.78%  "

$ws.Range("E25").Value = "  +0.10%  "

$ws.Range("D26").Value = "'71.01"
$ws.Range("E26").Value = "  -4.26%  "

$ws.Range("D27").Value = "'10.15"
$ws.Range("E27").Value = "  +2.32%  "

$ws.Range("D28").Value = "'2.816.48"
$ws.Range("E28").Value = "  -0.78%  "

$ws.Range("E29").Value = "  -3.15%  "

$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  -0.10%  "

$ws.Range("D31").Value = "'558.20"
$ws.Range("E31").Value = "  -6.54%  "

$ws.Range("E32").Value = "  -3.47%  "

$ws.Range("E33").Value = "  -3.75%  "

$ws.Range("D34").Value = "'1.94"
$ws.Range("E34").Value = "  -0.81%  "

$ws.Range("E35").Value = "  -1.79%  "

$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.08%  "

$ws.Range("E37").Value = "  -4.60%  "

$ws.Range("D38").Value = "'19.55"
$ws.Range("E38").Value = "  -1.72%  "

$ws.Range("D39").Value = "'155.80"
$ws.Range("E39").Value = "  -2.83%  "

$ws.Range("D40").Value = "'0.374"
$ws.Range("E40").Value = "  -1.80%  "

$ws.Range("E41").Value = "  -2.47%  "

$ws.Range("E42").Value = "  -4.49%  "

$ws.Range("D43").Value = "'17.95"
$ws.Range("E43").Value = "  -0.33%  "

$ws.Range("E44").Value = "  -6.81%  "

$ws.Range("E45").Value = "  +0.01%  "

$ws.Range("D46").Value = "'40.31"
$ws.Range("E46").Value = "  -1.01%  "

$ws.Range("E47").Value = "  -5.67%  "

$ws.Range("D48").Value = "'0.592"
$ws.Range("E48").Value = "  -2.41%  "

$ws.Range("D49").Value = "'153.74"
$ws.Range("E49").Value = "  -3.05%  "

$ws.Range("E50").Value = "  -2.21%  "

$ws.Range("E51").Value = "  -3.12%  "
